$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; existing rows 1..51 shift down to 3..53
$ws.Range("A1:A2").EntireRow.Insert()

# --- New row 1 ---
$ws.Range("A1").Value = 41725
$ws.Range("B1").Value = "  TRANSFERENCIA INTERNET"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "0000776328"
$ws.Range("E1").Value = "AG. NORTE"
$ws.Range("F1").Value = "380.88  "
$ws.Range("G1").Value = "2188.42"

# --- New row 2 ---
$ws.Range("A2").Value = 41725
$ws.Range("B2").Value = "RETIRO ATM BP D/KENNEDY 3"
$ws.Range("C2").Value = "D"
$ws.Range("D2").Value = "0000464551"
$ws.Range("E2").Value = "KENNEDY"
$ws.Range("F2").Value = "100.00  "
$ws.Range("G2").Value = "1807.54"

# Apply the same styles used for the other data rows in these columns
$ws.Range("A1:A2").Style = $ws.Range("A3").Style
$ws.Range("D1:D2").Style = $ws.Range("D3").Style
$ws.Range("F1:F2").Style = $ws.Range("F3").Style
$ws.Range("G1:G2").Style = $ws.Range("G3").Style

# Extend the shared formula in column H down to the (now pushed-down)
# rows that used to be rows 1 and 2.
$ws.Range("H1").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",TRIM(F1),", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL, ''mo_borrado_logico'' => false),")'
$ws.Range("H1:H4").FillDown()

# Two new blank rows are also present at the bottom of the sheet now
$ws.Range("A52").Value = $null
$ws.Range("A53").Value = $null
